# Weekly fruit/vegetable price update:
# Insert a new row at row 34 (pushing existing rows 34-40 down to 35-41)
# and populate it with the latest "Rabanito" price entry for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(34).Insert()

$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 44476
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 300000001
$ws.Cells.Item(34, 7).Value = "Rabanito"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 30
$ws.Cells.Item(34, 11).Value = 8000
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 8000
$ws.Cells.Item(34, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(34, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(34, 16).Value = 667
$ws.Cells.Item(34, 17).Value = 12
$ws.Cells.Item(34, 18).Value = "Hortaliza"
